$d = $word.ActiveDocument

# Locate the target paragraph: the last paragraph in the main body, which
# carries a custom left tab stop at 2110 twips (105.5 pt). It is currently
# empty and is the very last paragraph before the section break.
$n = $d.Paragraphs.Count
$target = $null
for ($i = $n; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $ts = $p.Format.TabStops
    $hasCustom = $false
    for ($j = 1; $j -le $ts.Count; $j++) {
        $pos = $ts.Item($j).Position
        if ([Math]::Abs($pos - 105.5) -lt 0.01) {
            $hasCustom = $true
            break
        }
    }
    if ($hasCustom) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph (tab stop at 2110)."
}

# Insert a fresh empty paragraph right before the target, then fill that new
# paragraph's range with the exact OOXML for the new "3 - SQL" section.
$prev = $target.Previous()
$insRange = $prev.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$newPara = $prev.Next()
$newRange = $newPara.Range

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>SQL</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1 –</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Qual o nome do retalhista (ou retalhistas) responsáveis pela reposição do maior número de categorias?</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">Nesta situação, considerámos que era pretendido o nome dos retalhistas com maior número de participações diretas na relação “responsavel_por”, não sendo </w:t></w:r><w:r><w:t xml:space="preserve">contabilizadas as </w:t></w:r><w:r><w:t>sub categorias das possíveis super categorias da sua responsabilidade.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>2 – Qual o nome do ou dos retalhistas que são responsáveis por todas as categorias simpels?</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:t>Neste caso, considerámos que para um retalhista ser responsável por todas as categorias simples tinha que</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>para cada categoria simples existente na base de dados, existir uma entrada na tabela “responsavel_por” que associasse essa categoria ao retalhista.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3 – Quais os produtos (ean) que nunca foram repostos?</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">Nesta </w:t></w:r><w:r><w:t>query</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>seleci</w:t></w:r><w:r><w:t>onámos todos os produtos cujos eans</w:t></w:r><w:r><w:t xml:space="preserve"> não constavam em nenhum evento de reposição.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xmlFrag)

# Add the vertAlign=superscript run-property default to the target
# paragraph's mark (its pPr/rPr), keeping the existing custom tab stop.
$target.Range.Font.Superscript = $true

Write-Host "Done"
